$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# "Yo" -> "Gabriel" (A7). Setting the value also drops the now-unused
# "Yo" shared string and appends "Gabriel" to the shared strings table.
$ws.Range("A7").Value = "Gabriel"

# Mark Prestacion/Convenio (B7:C7) with the same green font used by the
# Habitacion cell (C6) -- RGB 00B050.
$ws.Range("B7:C7").Font.Color = 5287936

# Move the active selection from C11 to C12.
$ws.Range("C12").Select()
